$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 12 new rows into the data table (it grows from 14 to 26
#    data rows) right before the old last-row / trailing blank rows,
#    pushing the "footer" block (signature lines) further down.
# ------------------------------------------------------------------
$ws.Rows("30:41").Insert()

# Row 41 is now the new "last" row of the table, so it should carry the
# special bottom-border formatting that used to live on row 29 (copy it
# out first, before row 29 itself gets restyled below).
$lastRow = $ws.Range("B29:J29")
$destLastRow = $ws.Range("B41:J41")
$lastRow.Copy()
$destLastRow.PasteSpecial(-4122)

# Copy the formatting of a "normal" data row (28) down into row 29 and
# the newly inserted rows 30-40 so they all pick up the same cell
# styles as the rest of the table body (row 29 used to be the special
# last row, but it no longer is).
$normalRow = $ws.Range("B28:J28")
for ($i = 29; $i -le 40; $i++) {
    $destRow = $ws.Range("B$i`:J$i")
    $normalRow.Copy()
    $destRow.PasteSpecial(-4122)
}

$ws.Application.CutCopyMode = $false

# ------------------------------------------------------------------
# 2. Update the summary header fields.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 1025141
$ws.Range("C13").Value = 6
$ws.Range("F13").Value = 23

# ------------------------------------------------------------------
# 3. Rewrite the worker/period detail table (rows 16-41).
# ------------------------------------------------------------------
$data = @(
    @("CC","1137220267","PEDRO JOSE TIRADO FERIA","2210",36341,908526),
    @("CC","1137220267","PEDRO JOSE TIRADO FERIA","2209",36341,908526),
    @("CC","73595113","MIGUEL ANGEL SUAREZ DE AGUSTIN","2210",36341,908526),
    @("CC","73595113","MIGUEL ANGEL SUAREZ DE AGUSTIN","2209",36341,908526),
    @("CC","72315354","JASIR ANTONIO DE LA HOZ OSPINO","1709",29509,1300000),
    @("CC","72315354","JASIR ANTONIO DE LA HOZ OSPINO","1708",29509,1300000),
    @("CC","72315354","JASIR ANTONIO DE LA HOZ OSPINO","1707",29509,1300000),
    @("CC","78115498","WILLIAN ANTONIO NORIEGA ROMERO","2002",44000,1100000),
    @("CC","92400566","MANUEL ANTONIO BARRAGAN PRIMERA","2011",33125,1100000),
    @("CC","92400566","MANUEL ANTONIO BARRAGAN PRIMERA","2010",33125,1100000),
    @("CC","92400566","MANUEL ANTONIO BARRAGAN PRIMERA","2009",33125,1100000),
    @("CC","92400566","MANUEL ANTONIO BARRAGAN PRIMERA","2008",33125,1100000),
    @("CC","92400566","MANUEL ANTONIO BARRAGAN PRIMERA","2007",33125,1100000),
    @("CC","92400566","MANUEL ANTONIO BARRAGAN PRIMERA","2006",33125,1100000),
    @("CC","92400566","MANUEL ANTONIO BARRAGAN PRIMERA","2005",33125,1100000),
    @("CC","92400566","MANUEL ANTONIO BARRAGAN PRIMERA","2004",33125,1100000),
    @("CC","92400566","MANUEL ANTONIO BARRAGAN PRIMERA","2003",33125,1100000),
    @("CC","92400566","MANUEL ANTONIO BARRAGAN PRIMERA","2002",33125,1100000),
    @("CC","1002410920","ROBERTO CARLOS ORTIZ PEREZ","2507",52000,1300000),
    @("CC","1002410920","ROBERTO CARLOS ORTIZ PEREZ","2506",52000,1300000),
    @("CC","1002410920","ROBERTO CARLOS ORTIZ PEREZ","2505",52000,1300000),
    @("CC","1002410920","ROBERTO CARLOS ORTIZ PEREZ","2504",52000,1300000),
    @("CC","1002410920","ROBERTO CARLOS ORTIZ PEREZ","2503",52000,1300000),
    @("CC","1002410920","ROBERTO CARLOS ORTIZ PEREZ","2502",52000,1300000),
    @("CC","1002410920","ROBERTO CARLOS ORTIZ PEREZ","2501",52000,1300000),
    @("CC","1002410920","ROBERTO CARLOS ORTIZ PEREZ","2412",52000,1300000)
)

$row = 16
foreach ($item in $data) {
    $ws.Range("B$row").Value = $item[0]
    $ws.Range("C$row").Value = $item[1]
    $ws.Range("D$row").Value = $item[2]
    $ws.Range("E$row").Value = $item[3]
    $ws.Range("F$row").Value = $item[4]
    $ws.Range("G$row").Value = $item[5]
    $row = $row + 1
}
